$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Insert a new row above row 31 (old row 31 and everything below shifts down by one)
$ws.Rows.Item(31).Insert()

# Fill in the new user's details in the freshly inserted row 31
$ws.Cells.Item(31, 1).Value = "Oleg_Babak"
$ws.Cells.Item(31, 2).Value = "Password1!"
$ws.Cells.Item(31, 4).Value = "CUSTOM_USER"
$ws.Cells.Item(31, 5).Value = "Smoke Test User"
$ws.Cells.Item(31, 6).Value = "N"

$ws.Range("B27").Select()
